$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Rename header labels (row 1)
$ws.Range("B1").Value = "description"
$ws.Range("C1").Value = "qty"
$ws.Range("D1").Value = "pack"
$ws.Range("A1").Value = "art"

# Remove the last product row (article 1647)
$ws.Rows("4:4").Delete()

# Mirror Excel's post-delete selection on the now-empty row 4
[void]$ws.Rows("4:4").Select()
